$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2579443333333333
$ws.Range("H2").Value = 0.773833
$ws.Range("I2").Value = 0.05467096027587241
$ws.Range("J2").Value = 0.0546709602758724
$ws.Range("M2").Value = 34.08558233333333
$ws.Range("N2").Value = 102.256747
$ws.Range("O2").Value = 0.3376420874206352
$ws.Range("P2").Value = 0.3376420874206352
$ws.Range("Q2").Value = 8.792182811250109
$ws.Range("R2").Value = 79.12964530125099
$ws.Range("S2").Value = 0.01845921714883618
$ws.Range("T2").Value = 0.01845921714883618
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2579443333333333
$ws.Range("H3").Value = 0.773833
$ws.Range("I3").Value = 0.05467096027587241
$ws.Range("J3").Value = 0.0546709602758724
$ws.Range("N3").Value = 90.767769
$ws.Range("O3").Value = 0.2997065709089496
$ws.Range("P3").Value = 0.2997065709089496
$ws.Range("Q3").Value = 7.804343887619666
$ws.Range("R3").Value = 70.23909498857699
$ws.Range("S3").Value = 0.01638524603258112
$ws.Range("T3").Value = 0.01638524603258112
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2579443333333333
$ws.Range("H4").Value = 0.773833
$ws.Range("I4").Value = 0.05467096027587241
$ws.Range("J4").Value = 0.0546709602758724
$ws.Range("M4").Value = 8.754337666666666
$ws.Range("N4").Value = 26.263013
$ws.Range("O4").Value = 0.08671797990283496
$ws.Range("P4").Value = 0.08671797990283496
$ws.Range("Q4").Value = 2.258131793203222
$ws.Range("R4").Value = 20.323186138829
$ws.Range("S4").Value = 0.004740955234471792
$ws.Range("T4").Value = 0.004740955234471791
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2579443333333333
$ws.Range("H5").Value = 0.773833
$ws.Range("I5").Value = 0.05467096027587241
$ws.Range("J5").Value = 0.0546709602758724
$ws.Range("M5").Value = 27.85597433333334
$ws.Range("N5").Value = 83.56792300000001
$ws.Range("O5").Value = 0.2759333617675802
$ws.Range("P5").Value = 0.2759333617675802
$ws.Range("Q5").Value = 7.185290728762111
$ws.Range("R5").Value = 64.66761655885901
$ws.Range("S5").Value = 0.0150855418599833
$ws.Range("T5").Value = 0.0150855418599833
$ws.Range("I6").Value = 0.4403936734732808
$ws.Range("J6").Value = 0.4403936734732807
$ws.Range("M6").Value = 34.08558233333333
$ws.Range("N6").Value = 102.256747
$ws.Range("O6").Value = 0.3376420874206352
$ws.Range("P6").Value = 0.3376420874206352
$ws.Range("Q6").Value = 70.82410234897388
$ws.Range("R6").Value = 637.416921140765
$ws.Range("S6").Value = 0.1486954391983601
$ws.Range("T6").Value = 0.1486954391983601
$ws.Range("I7").Value = 0.4403936734732808
$ws.Range("J7").Value = 0.4403936734732807
$ws.Range("N7").Value = 90.767769
$ws.Range("O7").Value = 0.2997065709089496
$ws.Range("P7").Value = 0.2997065709089496
$ws.Range("R7").Value = 565.8004342226551
$ws.Range("S7").Value = 0.1319888777266726
$ws.Range("T7").Value = 0.1319888777266726
$ws.Range("I8").Value = 0.4403936734732808
$ws.Range("J8").Value = 0.4403936734732807
$ws.Range("M8").Value = 8.754337666666666
$ws.Range("N8").Value = 26.263013
$ws.Range("O8").Value = 0.08671797990283496
$ws.Range("P8").Value = 0.08671797990283496
$ws.Range("Q8").Value = 18.19004002449278
$ws.Range("R8").Value = 163.710360220435
$ws.Range("S8").Value = 0.03819004972559163
$ws.Range("T8").Value = 0.03819004972559162
$ws.Range("I9").Value = 0.4403936734732808
$ws.Range("J9").Value = 0.4403936734732807
$ws.Range("M9").Value = 27.85597433333334
$ws.Range("N9").Value = 83.56792300000001
$ws.Range("O9").Value = 0.2759333617675802
$ws.Range("P9").Value = 0.2759333617675802
$ws.Range("Q9").Value = 57.8800255756539
$ws.Range("R9").Value = 520.920230180885
$ws.Range("S9").Value = 0.1215193068226564
$ws.Range("T9").Value = 0.1215193068226563
$ws.Range("G10").Value = 2.217259
$ws.Range("H10").Value = 6.651777
$ws.Range("I10").Value = 0.4699451123575263
$ws.Range("J10").Value = 0.4699451123575263
$ws.Range("M10").Value = 34.08558233333333
$ws.Range("N10").Value = 102.256747
$ws.Range("O10").Value = 0.3376420874206352
$ws.Range("P10").Value = 0.3376420874206352
$ws.Range("Q10").Value = 75.57656419882431
$ws.Range("R10").Value = 680.1890777894189
$ws.Range("S10").Value = 0.1586732487095201
$ws.Range("T10").Value = 0.1586732487095201
$ws.Range("G11").Value = 2.217259
$ws.Range("H11").Value = 6.651777
$ws.Range("I11").Value = 0.4699451123575263
$ws.Range("J11").Value = 0.4699451123575263
$ws.Range("N11").Value = 90.767769
$ws.Range("O11").Value = 0.2997065709089496
$ws.Range("P11").Value = 0.2997065709089496
$ws.Range("Q11").Value = 67.085217575057
$ws.Range("R11").Value = 603.766958175513
$ws.Range("S11").Value = 0.1408456381400953
$ws.Range("T11").Value = 0.1408456381400953
$ws.Range("G12").Value = 2.217259
$ws.Range("H12").Value = 6.651777
$ws.Range("I12").Value = 0.4699451123575263
$ws.Range("J12").Value = 0.4699451123575263
$ws.Range("M12").Value = 8.754337666666666
$ws.Range("N12").Value = 26.263013
$ws.Range("O12").Value = 0.08671797990283496
$ws.Range("P12").Value = 0.08671797990283496
$ws.Range("Q12").Value = 19.41063398045566
$ws.Range("R12").Value = 174.695705824101
$ws.Range("S12").Value = 0.04075269080885548
$ws.Range("T12").Value = 0.04075269080885548
$ws.Range("G13").Value = 2.217259
$ws.Range("H13").Value = 6.651777
$ws.Range("I13").Value = 0.4699451123575263
$ws.Range("J13").Value = 0.4699451123575263
$ws.Range("M13").Value = 27.85597433333334
$ws.Range("N13").Value = 83.56792300000001
$ws.Range("O13").Value = 0.2759333617675802
$ws.Range("P13").Value = 0.2759333617675802
$ws.Range("Q13").Value = 61.76390979435234
$ws.Range("R13").Value = 555.8751881491711
$ws.Range("S13").Value = 0.1296735346990554
$ws.Range("T13").Value = 0.1296735346990554
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1650883333333333
$ws.Range("H14").Value = 0.495265
$ws.Range("I14").Value = 0.03499025389332058
$ws.Range("J14").Value = 0.03499025389332058
$ws.Range("M14").Value = 34.08558233333333
$ws.Range("N14").Value = 102.256747
$ws.Range("O14").Value = 0.3376420874206352
$ws.Range("P14").Value = 0.3376420874206352
$ws.Range("Q14").Value = 5.62713197810611
$ws.Range("R14").Value = 50.644187802955
$ws.Range("S14").Value = 0.01181418236391877
$ws.Range("T14").Value = 0.01181418236391877
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1650883333333333
$ws.Range("H15").Value = 0.495265
$ws.Range("I15").Value = 0.03499025389332058
$ws.Range("J15").Value = 0.03499025389332058
$ws.Range("N15").Value = 90.767769
$ws.Range("O15").Value = 0.2997065709089496
$ws.Range("P15").Value = 0.2997065709089496
$ws.Range("Q15").Value = 4.994899901531666
$ws.Range("R15").Value = 44.954099113785
$ws.Range("S15").Value = 0.01048680900960064
$ws.Range("T15").Value = 0.01048680900960063
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1650883333333333
$ws.Range("H16").Value = 0.495265
$ws.Range("I16").Value = 0.03499025389332058
$ws.Range("J16").Value = 0.03499025389332058
$ws.Range("M16").Value = 8.754337666666666
$ws.Range("N16").Value = 26.263013
$ws.Range("O16").Value = 0.08671797990283496
$ws.Range("P16").Value = 0.08671797990283496
$ws.Range("Q16").Value = 1.445239014827222
$ws.Range("R16").Value = 13.007151133445
$ws.Range("S16").Value = 0.003034284133916067
$ws.Range("T16").Value = 0.003034284133916067
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1650883333333333
$ws.Range("H17").Value = 0.495265
$ws.Range("I17").Value = 0.03499025389332058
$ws.Range("J17").Value = 0.03499025389332058
$ws.Range("M17").Value = 27.85597433333334
$ws.Range("N17").Value = 83.56792300000001
$ws.Range("O17").Value = 0.2759333617675802
$ws.Range("P17").Value = 0.2759333617675802
$ws.Range("Q17").Value = 4.598696376066112
$ws.Range("R17").Value = 41.388267384595
$ws.Range("S17").Value = 0.009654978385885108
$ws.Range("T17").Value = 0.009654978385885107
